# Complete test for table_stats macro
#
# The workbook has a duplicate "Table_Stats" test row (row 56) whose macro
# column reads "table_stats_test1" - the same name used by another row.
# This removes that duplicate row (shifting the following rows up by one)
# and appends a completed/renamed test row at the end of the table:
#   Table_Stats | Test stats table | table_stats_test

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate "table_stats_test1" row - this shifts rows 57:75 up to 56:74.
$ws.Rows(56).Delete()

# Append the completed test as the new last row (75).
$ws.Range("A75").Value = "Table_Stats"
$ws.Range("B75").Value = "Test stats table"
$ws.Range("C75").Value = "table_stats_test"

# Match the author's final selection/scroll position.
$ws.Range("C75").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 49
$win.ScrollColumn = 1
